# Adds three new worksheets to the workbook:
#   - CypherOutput_Message  (copy of the existing "Message" sheet)
#   - StatOutput            (new stat-bar counts table)
#   - StatOutput_Message    (two "Message"-style blocks; second block uses
#                             the updated stat-bar Cypher query)

$wb = $excel.ActiveWorkbook

$msg = $wb.Worksheets.Item("Message")

# Pull the 10 "Message" rows (A1:A10) once so we can reuse them below.
$msgRows = @()
for ($r = 1; $r -le 10; $r++) {
    $msgRows += , $msg.Cells.Item($r, 1).Value2
}

$statCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Lip and oropharyngeal neoplasms malignant :: Melanoma-mucosa/maxilla']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# --- CypherOutput_Message : exact copy of Message --------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"

for ($r = 1; $r -le 10; $r++) {
    $cypherOutputMessage.Cells.Item($r, 1).Value = $msgRows[$r - 1]
}

# --- StatOutput : stat-bar counts table -------------------------------------
$statOutput = $wb.Worksheets.Add($null, $cypherOutputMessage)
$statOutput.Name = "StatOutput"

$statOutput.Cells.Item(1, 1).Value = "number_of_files"
$statOutput.Cells.Item(1, 2).Value = "number_of_sample"
$statOutput.Cells.Item(1, 3).Value = "number_of_cases"
$statOutput.Cells.Item(1, 4).Value = "number_of_study"

# The counts are written as text (matching the source export's shared-string
# cells), so force the "@" text format before assigning the numeric-looking
# strings -- otherwise they would be auto-detected as numbers.
$statOutput.Range("A2:D2").NumberFormat = "@"
$statOutput.Cells.Item(2, 1).Value = "76"
$statOutput.Cells.Item(2, 2).Value = "8"
$statOutput.Cells.Item(2, 3).Value = "4"
$statOutput.Cells.Item(2, 4).Value = "1"

# --- StatOutput_Message : two Message-style blocks --------------------------
$statOutputMessage = $wb.Worksheets.Add($null, $statOutput)
$statOutputMessage.Name = "StatOutput_Message"

for ($r = 1; $r -le 10; $r++) {
    $statOutputMessage.Cells.Item($r, 1).Value = $msgRows[$r - 1]
}

for ($r = 1; $r -le 10; $r++) {
    if ($r -eq 8) {
        $statOutputMessage.Cells.Item($r + 10, 1).Value = $statCypher
    } else {
        $statOutputMessage.Cells.Item($r + 10, 1).Value = $msgRows[$r - 1]
    }
}

# Adding sheets shifts the active tab to the last one created; restore the
# original active/selected sheet ("CypherOutput") to match the source file.
$wb.Worksheets.Item("CypherOutput").Activate()
